$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Linear Regression" results column (column E) along with its
# data. This shifts the remaining columns (F..L) left by one, so the
# "XGBoost" column becomes E, and the "*NOTE" column (L) becomes K -
# matching the target layout.
[void]$ws.Columns.Item(5).Delete()

# Update the selected cell as recorded in the saved view state.
[void]$ws.Range("G10").Select()
